# Update deterministic comparison values for July ruleset rerun.
# Affects Trace11, Trace17, and Trace22 sheets: columns R3 (PowellWYRelease)
# and T3 (Compact Point Volume) get refreshed model output values.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "Trace17"; R3 = 10069807.79058357;  T3 = 10352228.094583571 },
    @{ Sheet = "Trace22"; R3 = 8587570.2347199731; T3 = 8728198.2647199742 },
    @{ Sheet = "Trace11"; R3 = 8617022.5542262942; T3 = 8767620.2242262959 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range("R3").Value = $u.R3
    $ws.Range("T3").Value = $u.T3
}
